# Add a new textbox shape to slide 1 containing the repo link, placed just
# below the title block, matching the author's "Add files via upload" edit.
#
# Target (from the OOXML diff) is a new <p:sp> appended as the last shape in
# slide 1's shape tree:
#   <p:cNvPr id="6" name="TextBox 5">
#   <a:off x="119743" y="3320243"/>  <a:ext cx="7315200" cy="369332"/>
#   <a:noFill/>
#   <a:bodyPr wrap="square"><a:spAutoFit/></a:bodyPr>
#   single bold run: "https://github.com/danbui/law_model/tree/gemini-integration"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The presentation already has 4 top-level shapes (ids 1,2,3,4,9 in use, the
# highest being 9), but PowerPoint's "next shape id" counter in this deck is
# independent of that and would hand out id=5 for the very first shape we
# add. The real file's new shape has id=6, so we burn one id first by adding
# a throwaway textbox and immediately deleting it; the counter still moves
# forward, and the shape we actually keep then lands on id=6 / "TextBox 5",
# matching the target exactly.
$bump = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$bump.Delete()

# Position/size are specified in points for AddTextbox; the target EMU
# values (119743, 3320243, 7315200, 369332) divided by 914400 EMU-per-inch
# * 72 points-per-inch (i.e. /12700) give the point values below.
$left = 119743 / 12700
$top = 3320243 / 12700
$width = 7315200 / 12700
$height = 369332 / 12700

$shp = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$shp.Name = "TextBox 5"

# Plain rectangle textbox with no fill, auto-fit to the single line of text.
$shp.Fill.Visible = $false
$shp.TextFrame.WordWrap = $true
$shp.TextFrame.AutoSize = 1

$tr = $shp.TextFrame.TextRange
$tr.Text = "https://github.com/danbui/law_model/tree/gemini-integration"
$tr.Font.Bold = $true
